$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Clear()

$ws.Range("A1").Value = "Код города"
$ws.Range("B1").Value = "Название"
$ws.Range("C1").Value = "Регион"

$ws.Range("A2").Value = 383
$ws.Range("B2").Value = "Новосибирск"
$ws.Range("C2").Value = "Новосибирская область"

$ws.Range("A3").Value = 88552
$ws.Range("B3").Value = "Набережные Челны"
$ws.Range("C3").Value = "Татарстан"

$ws.Columns.Item(1).ColumnWidth = 11.166666666666666
$ws.Columns.Item(2).ColumnWidth = 21.833333333333336
$ws.Columns.Item(3).ColumnWidth = 24.5

$ws.Range("B6").Select() | Out-Null
